# Added Week 15 simulations
$wb = $excel.ActiveWorkbook

# OFF sheet - Home (row 2) stats updated
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 194
$wsOff.Range("C2").Value = 138
$wsOff.Range("D2").Value = 34
$wsOff.Range("E2").Value = 14

# DEF sheet - Home (row 2) stats updated
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 188
$wsDef.Range("C2").Value = 120
$wsDef.Range("D2").Value = 44
$wsDef.Range("E2").Value = 24
$wsDef.Range("F2").Value = 3
$wsDef.Range("G2").Value = 4
